# Auto-generated edit script: applies cyclic rotation of observation row data
# among rows {134,137,140}, {135,138,141}, {136,139,142}, {150,152,154,151,153}
# per the target diff (A,B,D,E,F,G,H,Q,R,U,W,AI columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 134
$ws.Range("A134").Value = 111973669
$ws.Range("B134").Value = 90660
$ws.Range("D134").Value = 'NT'
$ws.Range("E134").Value = 4362
$ws.Range("F134").Value = 'Blå taggsvamp'
$ws.Range("G134").Value = 'Hydnellum caeruleum'
$ws.Range("H134").Value = '(Hornem.) P.Karst.'
$ws.Range("Q134").Value = 437912.8137109271
$ws.Range("R134").Value = 6953242.433193879
$ws.Range("U134").Value = 'Härjedalen'
$ws.Range("W134").Value = 'Vemdalen'
$ws.Range("AI134").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 135
$ws.Range("A135").Value = 111973675
$ws.Range("B135").Value = 90660
$ws.Range("D135").Value = 'NT'
$ws.Range("E135").Value = 4362
$ws.Range("F135").Value = 'Blå taggsvamp'
$ws.Range("G135").Value = 'Hydnellum caeruleum'
$ws.Range("H135").Value = '(Hornem.) P.Karst.'
$ws.Range("Q135").Value = 438160.5186564626
$ws.Range("R135").Value = 6953262.568950667
$ws.Range("U135").Value = 'Härjedalen'
$ws.Range("W135").Value = 'Vemdalen'
$ws.Range("AI135").Value = 'äldre renbetad fattigristallskog med lavfläck på torr moränmark'

# Row 136
$ws.Range("A136").Value = 111973715
$ws.Range("B136").Value = 90678
$ws.Range("D136").Value = 'LC'
$ws.Range("E136").Value = 4366
$ws.Range("F136").Value = 'Skarp dropptaggsvamp'
$ws.Range("G136").Value = 'Hydnellum peckii'
$ws.Range("H136").Value = 'Banker'
$ws.Range("Q136").Value = 437962.8122493967
$ws.Range("R136").Value = 6953212.111986059
$ws.Range("U136").Value = 'Härjedalen'
$ws.Range("W136").Value = 'Vemdalen'
$ws.Range("AI136").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 137
$ws.Range("A137").Value = 111973671
$ws.Range("B137").Value = 88032
$ws.Range("D137").Value = 'VU'
$ws.Range("E137").Value = 6276
$ws.Range("F137").Value = 'Goliatmusseron'
$ws.Range("G137").Value = 'Tricholoma matsutake'
$ws.Range("H137").Value = '(S.Ito & S.Imai) Singer'
$ws.Range("Q137").Value = 438033.4411253001
$ws.Range("R137").Value = 6953252.100307667
$ws.Range("U137").Value = 'Härjedalen'
$ws.Range("W137").Value = 'Vemdalen'
$ws.Range("AI137").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 138
$ws.Range("A138").Value = 111973801
$ws.Range("B138").Value = 90660
$ws.Range("D138").Value = 'NT'
$ws.Range("E138").Value = 4362
$ws.Range("F138").Value = 'Blå taggsvamp'
$ws.Range("G138").Value = 'Hydnellum caeruleum'
$ws.Range("H138").Value = '(Hornem.) P.Karst.'
$ws.Range("Q138").Value = 437774.9523556355
$ws.Range("R138").Value = 6953074.713081508
$ws.Range("U138").Value = 'Härjedalen'
$ws.Range("W138").Value = 'Vemdalen'
$ws.Range("AI138").Value = 'äldre renbetad lingon- och lavtallskog på torr moränmark'

# Row 139
$ws.Range("A139").Value = 111973706
$ws.Range("B139").Value = 90652
$ws.Range("D139").Value = 'NT'
$ws.Range("E139").Value = 3100
$ws.Range("F139").Value = 'Talltaggsvamp'
$ws.Range("G139").Value = 'Bankera fuligineoalba'
$ws.Range("H139").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q139").Value = 438039.4596956634
$ws.Range("R139").Value = 6953083.047228135
$ws.Range("U139").Value = 'Härjedalen'
$ws.Range("W139").Value = 'Vemdalen'
$ws.Range("AI139").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 140
$ws.Range("A140").Value = 111973701
$ws.Range("B140").Value = 90682
$ws.Range("D140").Value = 'NT'
$ws.Range("E140").Value = 2059
$ws.Range("F140").Value = 'Skrovlig taggsvamp'
$ws.Range("G140").Value = 'Hydnellum scabrosum'
$ws.Range("H140").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q140").Value = 438134.3562496312
$ws.Range("R140").Value = 6953093.202189791
$ws.Range("U140").Value = 'Härjedalen'
$ws.Range("W140").Value = 'Vemdalen'
$ws.Range("AI140").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 141
$ws.Range("A141").Value = 111973730
$ws.Range("B141").Value = 90682
$ws.Range("D141").Value = 'NT'
$ws.Range("E141").Value = 2059
$ws.Range("F141").Value = 'Skrovlig taggsvamp'
$ws.Range("G141").Value = 'Hydnellum scabrosum'
$ws.Range("H141").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q141").Value = 437841.7707476557
$ws.Range("R141").Value = 6953302.531346441
$ws.Range("U141").Value = 'Härjedalen'
$ws.Range("W141").Value = 'Vemdalen'
$ws.Range("AI141").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 142
$ws.Range("A142").Value = 111973733
$ws.Range("B142").Value = 90660
$ws.Range("D142").Value = 'NT'
$ws.Range("E142").Value = 4362
$ws.Range("F142").Value = 'Blå taggsvamp'
$ws.Range("G142").Value = 'Hydnellum caeruleum'
$ws.Range("H142").Value = '(Hornem.) P.Karst.'
$ws.Range("Q142").Value = 437876.3203048867
$ws.Range("R142").Value = 6953355.130729643
$ws.Range("U142").Value = 'Härjedalen'
$ws.Range("W142").Value = 'Vemdalen'
$ws.Range("AI142").Value = 'äldre renbetad lingontallskog med lavfläckar på torr moränmark'

# Row 150
$ws.Range("A150").Value = 111973729
$ws.Range("B150").Value = 90652
$ws.Range("D150").Value = 'NT'
$ws.Range("E150").Value = 3100
$ws.Range("F150").Value = 'Talltaggsvamp'
$ws.Range("G150").Value = 'Bankera fuligineoalba'
$ws.Range("H150").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q150").Value = 437849.4557670459
$ws.Range("R150").Value = 6953246.83996963
$ws.Range("U150").Value = 'Härjedalen'
$ws.Range("W150").Value = 'Vemdalen'
$ws.Range("AI150").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 151
$ws.Range("A151").Value = 111973648
$ws.Range("B151").Value = 90682
$ws.Range("D151").Value = 'NT'
$ws.Range("E151").Value = 2059
$ws.Range("F151").Value = 'Skrovlig taggsvamp'
$ws.Range("G151").Value = 'Hydnellum scabrosum'
$ws.Range("H151").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q151").Value = 438637.1037644488
$ws.Range("R151").Value = 6953601.158154471
$ws.Range("U151").Value = 'Härjedalen'
$ws.Range("W151").Value = 'Vemdalen'
$ws.Range("AI151").Value = 'äldre fattigristallskog på torr moränmark'

# Row 152
$ws.Range("A152").Value = 111973663
$ws.Range("B152").Value = 90660
$ws.Range("D152").Value = 'NT'
$ws.Range("E152").Value = 4362
$ws.Range("F152").Value = 'Blå taggsvamp'
$ws.Range("G152").Value = 'Hydnellum caeruleum'
$ws.Range("H152").Value = '(Hornem.) P.Karst.'
$ws.Range("Q152").Value = 438169.8244046976
$ws.Range("R152").Value = 6953366.599956161
$ws.Range("U152").Value = 'Härjedalen'
$ws.Range("W152").Value = 'Vemdalen'
$ws.Range("AI152").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 153
$ws.Range("A153").Value = 111973699
$ws.Range("B153").Value = 90654
$ws.Range("D153").Value = 'VU'
$ws.Range("E153").Value = 149
$ws.Range("F153").Value = 'Tallgråticka'
$ws.Range("G153").Value = 'Boletopsis grisea'
$ws.Range("H153").Value = '(Peck) Bondartsev & Singer'
$ws.Range("Q153").Value = 438167.2097542446
$ws.Range("R153").Value = 6953104.522341937
$ws.Range("U153").Value = 'Härjedalen'
$ws.Range("W153").Value = 'Vemdalen'
$ws.Range("AI153").Value = 'äldre renbetad fattigris- och lavtallskog på torr moränmark'

# Row 154
$ws.Range("A154").Value = 111973736
$ws.Range("B154").Value = 90660
$ws.Range("D154").Value = 'NT'
$ws.Range("E154").Value = 4362
$ws.Range("F154").Value = 'Blå taggsvamp'
$ws.Range("G154").Value = 'Hydnellum caeruleum'
$ws.Range("H154").Value = '(Hornem.) P.Karst.'
$ws.Range("Q154").Value = 437842.6892572012
$ws.Range("R154").Value = 6953497.608170006
$ws.Range("U154").Value = 'Berg'
$ws.Range("W154").Value = 'Åsarne'
$ws.Range("AI154").Value = 'äldre renbetad lingontallskog med lavfläckar på torr moränmark'
